$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns C and D (record_atd, average_simulation_TD) for rows 2-11,
# and column C for row 12 (average of record_atd), per corrected relevance markers.

$ws.Range("C2").Value = 65
$ws.Range("D2").Value = 55

$ws.Range("C3").Value = 21
$ws.Range("D3").Value = 97

$ws.Range("C4").Value = 70
$ws.Range("D4").Value = 49

$ws.Range("C5").Value = 117
$ws.Range("D5").Value = 99.5

$ws.Range("C6").Value = 113
$ws.Range("D6").Value = 101

$ws.Range("C7").Value = 120
$ws.Range("D7").Value = 97

$ws.Range("C8").Value = 55
$ws.Range("D8").Value = 35.5

$ws.Range("C9").Value = 82
$ws.Range("D9").Value = 73

$ws.Range("C10").Value = 85
$ws.Range("D10").Value = 85

$ws.Range("C11").Value = 56
$ws.Range("D11").Value = 35.5

$ws.Range("C12").Value = 78.40000000000001
